# Updated symbol list on Thu Dec 22 14:25:55 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Price (column D) updates
$ws.Range("D2").Value = "242.79"
$ws.Range("D3").Value = "22.06"
$ws.Range("D4").Value = "5.457"
$ws.Range("D5").Value = "0.05754"
$ws.Range("D7").Value = "6.336"
$ws.Range("D8").Value = "0.8117"
$ws.Range("D9").Value = "0.8593"
$ws.Range("D10").Value = "0.1444"
$ws.Range("D11").Value = "0.07338"
$ws.Range("D12").Value = "0.03049"
$ws.Range("D13").Value = "0.03117"
$ws.Range("D14").Value = "0.09386"
$ws.Range("D15").Value = "3.938"
$ws.Range("D16").Value = "0.001592"
$ws.Range("D17").Value = "0.04838"
$ws.Range("D18").Value = "0.0005855"
$ws.Range("D19").Value = "0.006374"
$ws.Range("D20").Value = "0.004129"
$ws.Range("D21").Value = "0.0009997"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").Value = "3.718"
$ws.Range("D24").Value = "2.189"
$ws.Range("D25").Value = "0.3271"
$ws.Range("D26").Value = "0.1272"
$ws.Range("D27").Value = "0.0004003"
$ws.Range("D40").Value = "0.03856"
$ws.Range("D41").Value = "0.006719"
$ws.Range("D42").Value = "0.1066"
$ws.Range("D43").Value = "0.003203"
$ws.Range("D44").Value = "0.007141"
$ws.Range("D45").Value = "0.00005594"
$ws.Range("D47").Value = "0.3803"
$ws.Range("D48").Value = "0.1440"

# Volume/Data text (column E) updates
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
